# KAMAU FAMILY LIST - leader_assignment fix
#
# The G column ("KUNGU" family leaders list) had two erroneous / duplicate
# entries - "NANCY KUNGU" (G4) and "ANN KUNGU" (G7). Remove them and let the
# remaining names shift up so the list is contiguous again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the remaining G-column names up into place (G4 and G7 are dropped).
$ws.Range("G4").Value  = "SAM KUNGU"            # was G5
$ws.Range("G5").Value  = "JOEL KUNGU"           # was G6
$ws.Range("G6").Value  = "VICTOR KUNGU"         # was G8
$ws.Range("G7").Value  = "TERESIAH KUNGU"       # was G9
$ws.Range("G8").Value  = "ABIGAEL"              # was G10
$ws.Range("G9").Value  = "JAYDEN/CALEB/FELICIAH" # was G11

# Those last two rows carried the red "Times New Roman" font - keep that
# formatting on the cells it moved into.
$ws.Range("G8").Font.Color = 255
$ws.Range("G9").Font.Color = 255

# The two now-duplicated trailing rows are empty - remove them so the used
# range shrinks back down from K11 to K9.
$ws.Range("G10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Rows(10).Delete()
$ws.Rows(10).Delete()

# Leave the selection where the author ended up after the edit.
$ws.Range("G16").Select() | Out-Null
